$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 9 (anchor G9=5487)
$ws.Range("H9").Value = 76.166664
$ws.Range("I9").Value = 42.5
$ws.Range("J9").Value = 93
$ws.Range("K9").Value = 42.5
$ws.Range("L9").Value = 93
$ws.Range("M9").Value = 126.5
$ws.Range("N9").Value = -431

# Row 28 (anchor G28=27772)
$ws.Range("H28").Value = 1609.8667
$ws.Range("I28").Value = 1272.091
$ws.Range("K28").Value = 1272.091
$ws.Range("M28").Value = -787.0909999999999

# Row 33 (anchor G33=5512)
$ws.Range("H33").Value = 447.53333
$ws.Range("I33").Value = 541.1
$ws.Range("J33").Value = 260.4
$ws.Range("K33").Value = 541.1
$ws.Range("L33").Value = 260.4
$ws.Range("M33").Value = -312.1
$ws.Range("N33").Value = -718.4

# Row 62 (anchor G62=27781)
$ws.Range("H62").Value = 3920.6667
$ws.Range("J62").Value = 4287
$ws.Range("L62").Value = 4287
$ws.Range("N62").Value = -5535

# Row 65 (anchor G65=27781)
$ws.Range("H65").Value = 3920.6667
$ws.Range("J65").Value = 4287
$ws.Range("L65").Value = 21435
$ws.Range("N65").Value = -27675

# Row 74 (anchor G74=5507)
$ws.Range("H74").Value = 4849.7
$ws.Range("I74").Value = 4356.857
$ws.Range("K74").Value = 4356.857
$ws.Range("M74").Value = -3420.857

# Row 76 (anchor G76=12602)
$ws.Range("H76").Value = 2771.1538
$ws.Range("I76").Value = 2483
$ws.Range("J76").Value = 3107.3333
$ws.Range("K76").Value = 2483
$ws.Range("L76").Value = 3107.3333
$ws.Range("M76").Value = -2168
$ws.Range("N76").Value = -3737.3333

# Row 77 (anchor G77=5507)
$ws.Range("H77").Value = 4849.7
$ws.Range("I77").Value = 4356.857
$ws.Range("K77").Value = 21784.285
$ws.Range("M77").Value = -17104.285

# Row 79 (anchor G79=12602)
$ws.Range("H79").Value = 2771.1538
$ws.Range("I79").Value = 2483
$ws.Range("J79").Value = 3107.3333
$ws.Range("K79").Value = 2483
$ws.Range("L79").Value = 3107.3333
$ws.Range("M79").Value = -1391
$ws.Range("N79").Value = -5291.3333

# Row 96 (anchor G96=19894)
$ws.Range("H96").Value = 3490.8333
$ws.Range("J96").Value = 7208.375
$ws.Range("L96").Value = 21625.125
$ws.Range("N96").Value = -24371.125

# Row 98 (anchor G98=36237)
$ws.Range("H98").Value = 1047.75
$ws.Range("I98").Value = 1155.3529
$ws.Range("J98").Value = 438
$ws.Range("K98").Value = 1155.3529
$ws.Range("L98").Value = 438
$ws.Range("M98").Value = 342.6470999999999
$ws.Range("N98").Value = -3434

# Row 103 (anchor G103=19909)
$ws.Range("H103").Value = 764.6667
$ws.Range("I103").Value = 744.25
$ws.Range("K103").Value = 2232.75
$ws.Range("M103").Value = -1646.75

# Row 112 (anchor G112=27960)
$ws.Range("H112").Value = 3883.3076
$ws.Range("J112").Value = 3883.3076
$ws.Range("L112").Value = 11649.9228
$ws.Range("N112").Value = -13865.9228

# Row 122 (anchor G122=36237)
$ws.Range("H122").Value = 1047.75
$ws.Range("I122").Value = 1155.3529
$ws.Range("J122").Value = 438
$ws.Range("K122").Value = 3466.0587
$ws.Range("L122").Value = 1314
$ws.Range("M122").Value = -1016.0587
$ws.Range("N122").Value = -6214

# Row 132 (anchor G132=44049)
$ws.Range("H132").Value = 3136.2827
$ws.Range("I132").Value = 1066.3096
$ws.Range("K132").Value = 3198.9288
$ws.Range("M132").Value = -668.9288000000001

# Row 135 (anchor G135=44047)
$ws.Range("H135").Value = 2953.2354
$ws.Range("I135").Value = 954.9286
$ws.Range("J135").Value = 12278.667
$ws.Range("K135").Value = 8594.357399999999
$ws.Range("L135").Value = 110508.003
$ws.Range("M135").Value = -6059.357399999999
$ws.Range("N135").Value = -115578.003

# Row 137 (anchor G137=44013)
$ws.Range("H137").Value = 13208196
$ws.Range("I137").Value = 590749.44
$ws.Range("K137").Value = 1772248.32
$ws.Range("M137").Value = -1769698.32

# Row 138 (anchor G138=44169)
$ws.Range("H138").Value = 5480.827
$ws.Range("I138").Value = 1598.619
$ws.Range("J138").Value = 6839.6
$ws.Range("K138").Value = 4795.857
$ws.Range("L138").Value = 20518.8
$ws.Range("M138").Value = 344.143
$ws.Range("N138").Value = -30798.8

# Row 141 (anchor G141=44161)
$ws.Range("H141").Value = 2500
$ws.Range("I141").Value = 2500
$ws.Range("K141").Value = 7500
$ws.Range("M141").Value = -2320

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 4 (anchor G4=5071)
$ws.Range("H4").Value = 1367.3334
$ws.Range("J4").Value = 1367.3334
$ws.Range("L4").Value = 1367.3334
$ws.Range("N4").Value = -1599.3334

# Row 32 (anchor G32=44147)
$ws.Range("H32").Value = 20047.404
$ws.Range("I32").Value = 19666.564
$ws.Range("K32").Value = 19666.564
$ws.Range("M32").Value = -19379.564

# Row 61 (anchor G61=43999)
$ws.Range("H61").Value = 6045.5713
$ws.Range("I61").Value = 6619.8335
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 6619.8335
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -6407.8335
$ws.Range("N61").Value = -3024

# Row 122 (anchor G122=36168)
$ws.Range("H122").Value = 3788.425
$ws.Range("I122").Value = 3702.5
$ws.Range("J122").Value = 5421
$ws.Range("K122").Value = 11107.5
$ws.Range("L122").Value = 16263
$ws.Range("M122").Value = -8657.5
$ws.Range("N122").Value = -21163

# Row 132 (anchor G132=43997)
$ws.Range("H132").Value = 1957.4615
$ws.Range("I132").Value = 1873.3889
$ws.Range("K132").Value = 5620.1667
$ws.Range("M132").Value = -3090.1667

# Row 136 (anchor G136=43999)
$ws.Range("H136").Value = 6045.5713
$ws.Range("I136").Value = 6619.8335
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 19859.5005
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -17309.5005
$ws.Range("N136").Value = -12900

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 94 (anchor G94=19939)
$ws.Range("H94").Value = 1090.8334
$ws.Range("I94").Value = 909.4
$ws.Range("J94").Value = 1998
$ws.Range("K94").Value = 909.4
$ws.Range("L94").Value = 1998
$ws.Range("M94").Value = -458.4
$ws.Range("N94").Value = -2900

# Row 105 (anchor G105=19947)
$ws.Range("H105").Value = 4877.6875
$ws.Range("I105").Value = 4723.857
$ws.Range("K105").Value = 4723.857
$ws.Range("M105").Value = -2976.857

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 16 (anchor G16=27691)
$ws.Range("H16").Value = 896.0769
$ws.Range("I16").Value = 785.5833
$ws.Range("K16").Value = 785.5833
$ws.Range("M16").Value = -498.5833

# Row 22 (anchor G22=5367)
$ws.Range("H22").Value = 1767.7142
$ws.Range("I22").Value = 2011.5
$ws.Range("J22").Value = 1670.2
$ws.Range("K22").Value = 2011.5
$ws.Range("L22").Value = 1670.2
$ws.Range("M22").Value = -1661.5
$ws.Range("N22").Value = -2370.2

# Row 31 (anchor G31=44023)
$ws.Range("H31").Value = 15627558
$ws.Range("I31").Value = 18183866
$ws.Range("J31").Value = 5674.8887
$ws.Range("K31").Value = 18183866
$ws.Range("L31").Value = 5674.8887
$ws.Range("M31").Value = -18183571
$ws.Range("N31").Value = -6264.8887

# Row 34 (anchor G34=44023)
$ws.Range("H34").Value = 15627558
$ws.Range("I34").Value = 18183866
$ws.Range("J34").Value = 5674.8887
$ws.Range("K34").Value = 18183866
$ws.Range("L34").Value = 5674.8887
$ws.Range("M34").Value = -18183664
$ws.Range("N34").Value = -6078.8887

# Row 60 (anchor G60=1937)
$ws.Range("H60").Value = 30000
$ws.Range("J60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -31022

# Row 94 (anchor G94=32934)
$ws.Range("H94").Value = 2280.4666
$ws.Range("I94").Value = 1073.5
$ws.Range("J94").Value = 2719.3635
$ws.Range("K94").Value = 1073.5
$ws.Range("L94").Value = 2719.3635
$ws.Range("M94").Value = -622.5
$ws.Range("N94").Value = -3621.3635

# Row 99 (anchor G99=36198)
$ws.Range("H99").Value = 13764.4
$ws.Range("J99").Value = 7083
$ws.Range("L99").Value = 7083
$ws.Range("N99").Value = -10079

# Row 107 (anchor G107=27689)
$ws.Range("H107").Value = 941.5238
$ws.Range("I107").Value = 579.2857
$ws.Range("K107").Value = 579.2857
$ws.Range("M107").Value = 1340.7143

# Row 113 (anchor G113=27691)
$ws.Range("H113").Value = 896.0769
$ws.Range("I113").Value = 785.5833
$ws.Range("K113").Value = 785.5833
$ws.Range("M113").Value = 1384.4167

# Row 126 (anchor G126=36198)
$ws.Range("H126").Value = 13764.4
$ws.Range("J126").Value = 7083
$ws.Range("L126").Value = 21249
$ws.Range("N126").Value = -26189

# Row 132 (anchor G132=44019)
$ws.Range("H132").Value = 102573496
$ws.Range("I132").Value = 111112950
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 333338850
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -333336320
$ws.Range("N132").Value = -305060

# Row 134 (anchor G134=44020)
$ws.Range("H134").Value = 2556.389
$ws.Range("I134").Value = 1705.8
$ws.Range("K134").Value = 5117.4
$ws.Range("M134").Value = -2582.4

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 2 (anchor G2=4847)
$ws.Range("H2").Value = 43.875
$ws.Range("J2").Value = 53.75
$ws.Range("L2").Value = 322.5
$ws.Range("N2").Value = -548.5

# Row 12 (anchor G12=4854)
$ws.Range("H12").Value = 193.13333
$ws.Range("I12").Value = 246.6
$ws.Range("K12").Value = 739.8
$ws.Range("M12").Value = -566.8

# Row 37 (anchor G37=9516)
$ws.Range("H37").Value = 100106980
$ws.Range("J37").Value = 100106980
$ws.Range("L37").Value = 300320940
$ws.Range("N37").Value = -300321164

# Row 55 (anchor G55=4733)
$ws.Range("H55").Value = 2800.8
$ws.Range("J55").Value = 3333.3333
$ws.Range("L55").Value = 9999.999899999999
$ws.Range("N55").Value = -10353.9999

# Row 114 (anchor G114=27865)
$ws.Range("H114").Value = 1868.1428
$ws.Range("I114").Value = 125
$ws.Range("K114").Value = 375
$ws.Range("M114").Value = 2879

# Row 129 (anchor G129=36054)
$ws.Range("H129").Value = 1844.3462
$ws.Range("I129").Value = 1450.6666
$ws.Range("J129").Value = 2052.7646
$ws.Range("K129").Value = 4351.9998
$ws.Range("L129").Value = 6158.293799999999
$ws.Range("M129").Value = 648.0002000000004
$ws.Range("N129").Value = -16158.2938

# Row 131 (anchor G131=36060)
$ws.Range("H131").Value = 5145554.5
$ws.Range("J131").Value = 6174743.5
$ws.Range("L131").Value = 18524230.5
$ws.Range("N131").Value = -18534310.5

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 97 (anchor G97=19940)
$ws.Range("H97").Value = 2078.647
$ws.Range("I97").Value = 1455.5555
$ws.Range("K97").Value = 1455.5555
$ws.Range("M97").Value = -959.5554999999999

# Row 102 (anchor G102=36169)
$ws.Range("H102").Value = 11115079
$ws.Range("I102").Value = 13517682
$ws.Range("K102").Value = 13517682
$ws.Range("M102").Value = -13516060

# Row 122 (anchor G122=36182)
$ws.Range("H122").Value = 230587.02
$ws.Range("I122").Value = 371819.38
$ws.Range("K122").Value = 1115458.14
$ws.Range("M122").Value = -1113008.14

# Row 132 (anchor G132=44008)
$ws.Range("H132").Value = 136194.67
$ws.Range("I132").Value = 201262
$ws.Range("K132").Value = 603786
$ws.Range("M132").Value = -601256

# Row 139 (anchor G139=42373)
$ws.Range("H139").Value = 88362.5
$ws.Range("J139").Value = 88362.5
$ws.Range("L139").Value = 88362.5
$ws.Range("N139").Value = -98642.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 7 (anchor G7=36249)
$ws.Range("H7").Value = 4669.558
$ws.Range("I7").Value = 3861.4285
$ws.Range("K7").Value = 3861.4285
$ws.Range("M7").Value = -3749.4285

# Row 22 (anchor G22=5277)
$ws.Range("H22").Value = 700.94116
$ws.Range("I22").Value = 819.3333
$ws.Range("J22").Value = 636.36365
$ws.Range("K22").Value = 819.3333
$ws.Range("L22").Value = 636.36365
$ws.Range("M22").Value = -524.3333
$ws.Range("N22").Value = -1226.36365

# Row 27 (anchor G27=5277)
$ws.Range("H27").Value = 700.94116
$ws.Range("I27").Value = 819.3333
$ws.Range("J27").Value = 636.36365
$ws.Range("K27").Value = 819.3333
$ws.Range("L27").Value = 636.36365
$ws.Range("M27").Value = -712.3333
$ws.Range("N27").Value = -850.36365

# Row 40 (anchor G40=36248)
$ws.Range("H40").Value = 16371881
$ws.Range("I40").Value = 5954183
$ws.Range("J40").Value = 47624976
$ws.Range("K40").Value = 5954183
$ws.Range("L40").Value = 47624976
$ws.Range("M40").Value = -5954047
$ws.Range("N40").Value = -47625248

# Row 46 (anchor G46=5282)
$ws.Range("H46").Value = 4991.4287
$ws.Range("I46").Value = 3001
$ws.Range("K46").Value = 3001
$ws.Range("M46").Value = -2813

# Row 55 (anchor G55=5284)
$ws.Range("H55").Value = 353.25
$ws.Range("I55").Value = 122.888885
$ws.Range("J55").Value = 649.4286
$ws.Range("K55").Value = 122.888885
$ws.Range("L55").Value = 649.4286
$ws.Range("M55").Value = 50.111115
$ws.Range("N55").Value = -995.4286

# Row 61 (anchor G61=27740)
$ws.Range("H61").Value = 11401
$ws.Range("I61").Value = 14899
$ws.Range("K61").Value = 14899
$ws.Range("M61").Value = -14697

# Row 68 (anchor G68=12563)
$ws.Range("H68").Value = 4750.1816
$ws.Range("I68").Value = 3795
$ws.Range("J68").Value = 7297.3335
$ws.Range("K68").Value = 3795
$ws.Range("L68").Value = 7297.3335
$ws.Range("M68").Value = -3046
$ws.Range("N68").Value = -8795.3335

# Row 71 (anchor G71=12563)
$ws.Range("H71").Value = 4750.1816
$ws.Range("I71").Value = 3795
$ws.Range("J71").Value = 7297.3335
$ws.Range("K71").Value = 18975
$ws.Range("L71").Value = 36486.6675
$ws.Range("M71").Value = -15231
$ws.Range("N71").Value = -43974.6675

# Row 113 (anchor G113=27740)
$ws.Range("H113").Value = 11401
$ws.Range("I113").Value = 14899
$ws.Range("K113").Value = 14899
$ws.Range("M113").Value = -12729

# Row 122 (anchor G122=36247)
$ws.Range("H122").Value = 7740.8066
$ws.Range("I122").Value = 4447.1763
$ws.Range("J122").Value = 11740.214
$ws.Range("K122").Value = 13341.5289
$ws.Range("L122").Value = 35220.642
$ws.Range("M122").Value = -10891.5289
$ws.Range("N122").Value = -40120.642

# Row 126 (anchor G126=36249)
$ws.Range("H126").Value = 4669.558
$ws.Range("I126").Value = 3861.4285
$ws.Range("K126").Value = 11584.2855
$ws.Range("M126").Value = -9114.2855

# Row 132 (anchor G132=44058)
$ws.Range("H132").Value = 2764.82
$ws.Range("I132").Value = 2746.9443
$ws.Range("J132").Value = 2810.7856
$ws.Range("K132").Value = 8240.832900000001
$ws.Range("L132").Value = 8432.356800000001
$ws.Range("M132").Value = -5710.832900000001
$ws.Range("N132").Value = -13492.3568

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 107 (anchor G107=27746)
$ws.Range("H107").Value = 738.41174
$ws.Range("I107").Value = 572.55554
$ws.Range("J107").Value = 925
$ws.Range("K107").Value = 1717.66662
$ws.Range("L107").Value = 2775
$ws.Range("M107").Value = 202.33338
$ws.Range("N107").Value = -6615

# Row 122 (anchor G122=36208)
$ws.Range("H122").Value = 2546.7917
$ws.Range("I122").Value = 885.75
$ws.Range("J122").Value = 4207.8335
$ws.Range("K122").Value = 2657.25
$ws.Range("L122").Value = 12623.5005
$ws.Range("M122").Value = -207.25
$ws.Range("N122").Value = -17523.5005

# Row 126 (anchor G126=36210)
$ws.Range("H126").Value = 2024.1177
$ws.Range("I126").Value = 1964.5555
$ws.Range("K126").Value = 5893.666499999999
$ws.Range("M126").Value = -3423.666499999999

# Row 132 (anchor G132=44029)
$ws.Range("H132").Value = 333333340
$ws.Range("I132").Value = 333333340
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1000000020
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -999997490
$ws.Range("N132").ClearContents()

# Row 136 (anchor G136=44031)
$ws.Range("H136").Value = 2953.6956
$ws.Range("J136").Value = 6887.778
$ws.Range("L136").Value = 20663.334
$ws.Range("N136").Value = -25763.334
